$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.026.12"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "3.500.71"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'604.99"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "'172.48"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D8").Value = "3.495.77"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("D11").Value = "'7.21"
$ws.Range("E11").Value = "  +6.77%  "
$ws.Range("D12").Value = "'0.586"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'46.01"
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("D14").Value = "'0.0000275"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").Value = "4.069.06"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "'8.35"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "'612.48"
$ws.Range("D18").Value = "3.512.26"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "70.027.68"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").Value = "'17.50"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "'0.877"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "'9.14"
$ws.Range("E23").Value = "  -8.11%  "
$ws.Range("D24").Value = "'98.58"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "'15.50"
$ws.Range("E25").Value = "  -2.82%  "
$ws.Range("D26").Value = "'3.71"
$ws.Range("E26").Value = "  -3.56%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "'2.55"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").Value = "'33.74"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").Value = "'8.98"
$ws.Range("E30").Value = "  -2.77%  "
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").Value = "'2.98"
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'8.04"
$ws.Range("E32").Value = "  -5.13%  "
$ws.Range("D33").Value = "'1.27"
$ws.Range("E33").Value = "  -4.92%  "
$ws.Range("D34").Value = "'630.99"
$ws.Range("E34").Value = "  +11.09%  "
$ws.Range("D35").Value = "'6.79"
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("D36").Value = "'0.0996"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").Value = "'10.72"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("E38").Value = "  +5.45%  "
$ws.Range("D39").Value = "'3.46"
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("D40").Value = "'56.73"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "3.357.40"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("D45").Value = "'0.309"
$ws.Range("E45").Value = "  -5.84%  "
$ws.Range("D46").Value = "'2.90"
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").Value = "'31.81"
$ws.Range("E47").Value = "  -3.87%  "
$ws.Range("E48").Value = "  -4.30%  "
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").Value = "'132.96"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("E51").Value = "  -0.02%  "
